$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K3").Value = "2016-08-25 06:47:13"
$wsZh.Range("P3").Value = ""
$wsZh.Columns.Item(16).AutoFit()

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K3").Value = "2016-08-25 06:47:20"
$wsDe.Range("P3").Value = ""
$wsDe.Columns.Item(16).AutoFit()
